{"js": "// Update the three-digit-divided-by-one-digit practice table:\n// replace each division-problem cell's text with its new prompt.\n// Cells are addressed positionally (row, col) in document order so\n// that the edit is correct even though some new values coincide with\n// other (unrelated) old values elsewhere in the table.\nconst replacements = [\n  \"472\u00f73=\", \"382\u00f72=\",\n  \"780\u00f73=\", \"529\u00f78=\",\n  \"767\u00f77=\", \"704\u00f72=\",\n  \"998\u00f75=\", \"424\u00f78=\",\n  \"169\u00f76=\", \"145\u00f76=\",\n  \"440\u00f78=\", \"295\u00f73=\",\n  \"821\u00f79=\", \"332\u00f79=\",\n  \"362\u00f77=\", \"381\u00f75=\",\n  \"516\u00f79=\", \"606\u00f77=\",\n  \"560\u00f73=\", \"660\u00f76=\",\n  \"861\u00f77=\", \"885\u00f74=\",\n  \"820\u00f77=\", \"594\u00f72=\",\n  \"551\u00f79=\", \"766\u00f76=\",\n  \"370\u00f75=\", \"949\u00f75=\",\n  \"672\u00f73=\", \"271\u00f77=\",\n  \"897\u00f77=\", \"729\u00f76=\",\n  \"236\u00f73=\", \"710\u00f73=\",\n  \"885\u00f74=\", \"227\u00f76=\",\n  \"595\u00f72=\", \"780\u00f72=\",\n  \"664\u00f73=\", \"773\u00f78=\",\n  \"321\u00f76=\", \"947\u00f73=\",\n  \"810\u00f78=\", \"820\u00f79=\",\n  \"418\u00f76=\", \"215\u00f72=\",\n  \"207\u00f75=\", \"313\u00f76=\",\n  \"838\u00f72=\", \"918\u00f79=\",\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Gather every cell's current text in document order (row-major),\n// loading all of them in one batch.\nconst allCells = [];\nfor (const row of rows.items) {\n  const cells = row.cells;\n  cells.load(\"items\");\n  allCells.push(cells);\n}\nawait context.sync();\n\nconst cellList = [];\nfor (const cells of allCells) {\n  for (const cell of cells.items) {\n    cell.load(\"value\");\n    cellList.push(cell);\n  }\n}\nawait context.sync();\n\n// Walk the cells in order and consume the replacement list one entry\n// at a time (so a new value that matches a later old value is never\n// mistaken for a fresh match).\nlet pos = 0;\nfor (const cell of cellList) {\n  if (pos >= replacements.length) break;\n  const expectedOld = replacements[pos];\n  if (cell.value === expectedOld) {\n    cell.value = replacements[pos + 1];\n    pos += 2;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the three-digit-divided-by-one-digit practice table:\n# replace each division-problem cell's text with its new prompt.\n# Cells are addressed positionally (row, col) in document order so\n# that the edit is correct even though some new values coincide with\n# other (unrelated) old values elsewhere in the table.\n$replacements = @(\n  \"472\u00f73=\", \"382\u00f72=\",\n  \"780\u00f73=\", \"529\u00f78=\",\n  \"767\u00f77=\", \"704\u00f72=\",\n  \"998\u00f75=\", \"424\u00f78=\",\n  \"169\u00f76=\", \"145\u00f76=\",\n  \"440\u00f78=\", \"295\u00f73=\",\n  \"821\u00f79=\", \"332\u00f79=\",\n  \"362\u00f77=\", \"381\u00f75=\",\n  \"516\u00f79=\", \"606\u00f77=\",\n  \"560\u00f73=\", \"660\u00f76=\",\n  \"861\u00f77=\", \"885\u00f74=\",\n  \"820\u00f77=\", \"594\u00f72=\",\n  \"551\u00f79=\", \"766\u00f76=\",\n  \"370\u00f75=\", \"949\u00f75=\",\n  \"672\u00f73=\", \"271\u00f77=\",\n  \"897\u00f77=\", \"729\u00f76=\",\n  \"236\u00f73=\", \"710\u00f73=\",\n  \"885\u00f74=\", \"227\u00f76=\",\n  \"595\u00f72=\", \"780\u00f72=\",\n  \"664\u00f73=\", \"773\u00f78=\",\n  \"321\u00f76=\", \"947\u00f73=\",\n  \"810\u00f78=\", \"820\u00f79=\",\n  \"418\u00f76=\", \"215\u00f72=\",\n  \"207\u00f75=\", \"313\u00f76=\",\n  \"838\u00f72=\", \"918\u00f79=\"\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n$rows = $t.Rows.Count\n$cols = $t.Columns.Count\n\n$pos = 0\nfor ($r = 1; $r -le $rows; $r++) {\n  for ($c = 1; $c -le $cols; $c++) {\n    if ($pos -ge $replacements.Length) { continue }\n    $cell = $t.Cell($r, $c)\n    $raw = $cell.Range.Text\n    $clean = $raw.TrimEnd([char]13, [char]7)\n    $expectedOld = $replacements[$pos]\n    if ($clean -eq $expectedOld) {\n      $cell.Range.Text = $replacements[$pos + 1]\n      $pos += 2\n    }\n  }\n}\n"}
